# Add new expense row (row 27): Tomáš, "3x ESP32", 23.85, "techfun"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = "Tomáš"
$ws.Range("B27").Value = "3x ESP32"
$ws.Range("C27").Value = 23.85
$ws.Range("D27").Value = "techfun"

$ws.Range("C27").NumberFormat = "#,##0.00\ [$€-41B]"

$ws.Range("E27").Select()
